$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: new "Sprint Meeting" entries added to the two right-hand
#     mini-tables (M:Q "Tabelle24" and W:AA "Tabelle245") ---

# Tabelle24 (M:Q) - Aufgabe/Datum/Von/Bis/Status
$ws.Range("M19").Value = "Sprint Meeting"
$ws.Range("N19").Value = 42697
$ws.Range("O19").Value = 0.59027777777777779
$ws.Range("P19").Value = 0.65972222222222221
$ws.Range("Q19").Value = 1
$ws.Range("Q19").NumberFormat = "0%"

# Tabelle245 (W:AA) - Aufgabe/Datum/Von/Bis/Status
$ws.Range("W19").Value = "Sprint Meeting"
$ws.Range("X19").Value = 42697
$ws.Range("Y19").Value = 0.59027777777777779
$ws.Range("Z19").Value = 0.65972222222222221
$ws.Range("AA19").Value = 1
$ws.Range("AA19").NumberFormat = "0%"

# --- Row 22: new "Sprint Meeting" entry added to the left-hand
#     Tabelle2 (C:G) ---
$ws.Range("C22").Value = "Sprint Meeting"
$ws.Range("D22").Value = 42697
$ws.Range("E22").Value = 0.59027777777777779
$ws.Range("F22").Value = 0.65972222222222221
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = "0%"

# --- View state: move the window so row 10 is at the top and select
#     the newly added row in the W:AA table ---
$excel.Goto($ws.Range("A10"), $true)
$ws.Range("W19:AA19").Select()
